$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.13429999999998
$ws.Range("A7").Value = -20.21459999999999
$ws.Range("D7").Value = -7.014799999999997
$ws.Range("D15").Value = -8.607999999999997
$ws.Range("A16").Value = -21.77980000000001
$ws.Range("E16").Value = 16.06330000000001
$ws.Range("E19").Value = 16.3827
$ws.Range("D21").Value = -8.748799999999999
$ws.Range("D22").Value = -8.0402
$ws.Range("D23").Value = -7.119699999999995
$ws.Range("A28").Value = -21.98719999999999
$ws.Range("A29").Value = -21.27099999999997
$ws.Range("A32").Value = -21.21319999999999
$ws.Range("D34").Value = -7.9331
$ws.Range("E36").Value = 15.73610000000001
$ws.Range("A40").Value = -19.9962
$ws.Range("D43").Value = -8.396299999999998
$ws.Range("D45").Value = -7.725
$ws.Range("E46").Value = 17.10349999999999
$ws.Range("D50").Value = -8.4527
$ws.Range("E50").Value = 16.52660000000001
$ws.Range("D51").Value = -7.342199999999998
$ws.Range("A52").Value = -22.1949
$ws.Range("A57").Value = -22.3051
$ws.Range("A66").Value = -21.33379999999999
$ws.Range("D66").Value = -7.808400000000006
$ws.Range("D67").Value = -6.445700000000003
$ws.Range("D79").Value = -6.298500000000003
$ws.Range("D84").Value = -8.808100000000003
$ws.Range("D92").Value = -6.565600000000004
$ws.Range("E95").Value = 18.16400000000002
$ws.Range("D97").Value = -8.579500000000001
$ws.Range("E97").Value = 16.46419999999998
$ws.Range("A100").Value = -21.83610000000001
